$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 61873
$ws.Range("B2").Value = "Pietra Correia"
$ws.Range("C2").Value = "Financeiro"
$ws.Range("D2").Value = "Consulta medica"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45078
$ws.Range("G2").Value = 8308.49

# Row 3
$ws.Range("A3").Value = 79581
$ws.Range("B3").Value = "Sr. Luiz Felipe da Conceição"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Consulta medica"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45087
$ws.Range("G3").Value = 3749.39

# Row 4
$ws.Range("A4").Value = 4004
$ws.Range("B4").Value = "Mariana Duarte"
$ws.Range("C4").Value = "TI"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45083
$ws.Range("G4").Value = 3477.25

# Row 5
$ws.Range("A5").Value = 22174
$ws.Range("B5").Value = "Sr. João Pedro da Conceição"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Outros"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45092
$ws.Range("G5").Value = 9609.809999999999

# Row 6
$ws.Range("A6").Value = 32826
$ws.Range("B6").Value = "Davi Miguel Novais"
$ws.Range("C6").Value = "Engenharia"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 4025.66

# Row 7
$ws.Range("A7").Value = 96440
$ws.Range("B7").Value = "Guilherme Cassiano"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Doenca"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45087
$ws.Range("G7").Value = 5880.94

# Row 8
$ws.Range("A8").Value = 29404
$ws.Range("B8").Value = "João Pedro Camargo"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45088
$ws.Range("G8").Value = 4247.58

# Row 9
$ws.Range("A9").Value = 94305
$ws.Range("B9").Value = "Lucas Melo"
$ws.Range("C9").Value = "Operacoes"
$ws.Range("D9").Value = "Consulta medica"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45103
$ws.Range("G9").Value = 3628.86

# Row 10
$ws.Range("A10").Value = 28221
$ws.Range("B10").Value = "Thomas Melo"
$ws.Range("C10").Value = "Marketing"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 45106
$ws.Range("G10").Value = 3161

# Row 11
$ws.Range("A11").Value = 30396
$ws.Range("B11").Value = "Maria Vitória Abreu"
$ws.Range("C11").Value = "Atendimento ao Cliente"
$ws.Range("D11").Value = "Viagem de negocios"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45079
$ws.Range("G11").Value = 9695.41
